$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.130952380952381, 0.469047619047619, 0.5435714285714285, 0.7372729674763808, 0.7343449363085427, 42)
    3  = @(0.006862445399837091, 0.6039302007301544, 0.7234397629857844, 0.8505526221144606, 0.8610908701029418, 41)
    4  = @(0.1524205693483774, 0.6423595844912023, 0.8775854164232644, 0.9367952905642003, 0.9360875405385815, 40)
    5  = @(0.04980756042726093, 0.7196857601028602, 0.9118160530585012, 0.9548905974290988, 0.9660564865707858, 39)
    6  = @(0.1685415784894023, 0.7315745272820626, 0.9669495865848617, 0.983335947977527, 0.9817888564550711, 38)
    7  = @(0.05236864345981531, 0.6876603977094401, 0.8135759198832745, 0.9019844343907907, 0.9128836707474132, 37)
    8  = @(0.1715962988235919, 0.6741482612662726, 0.9067807233488873, 0.9522503469933142, 0.9499485040004092, 36)
    9  = @(0.1127370711831118, 0.5725331809908775, 0.6431051428235416, 0.801938365975554, 0.805565974699134, 35)
    10 = @(0.1706140055022492, 0.666628778494778, 0.750579926290491, 0.8663601596856189, 0.8621679294376572, 34)
    11 = @(0.2039437044841572, 0.6558282946641647, 0.7409461765722168, 0.8607823049832152, 0.8492396173402222, 33)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("B$row").Value = $values[0]
    $ws.Range("C$row").Value = $values[1]
    $ws.Range("D$row").Value = $values[2]
    $ws.Range("E$row").Value = $values[3]
    $ws.Range("F$row").Value = $values[4]
    $ws.Range("G$row").Value = $values[5]
}
